$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet0")

# Add new plate-loading notes for WEEK 3 (D14:D16)
$ws.Range("D14").Value = "4.5x1, 2.3x2"
$ws.Range("D16").Value = "11.3x2, 2.3x1"
$ws.Range("D15").Value = "11.3x1"

# Update existing D2:D4 values (WEEK 1 plate-loading notes)
$ws.Range("D2").Value = "1x4.5, 1x2.3"
$ws.Range("D3").Value = "2x4.5"
$ws.Range("D4").Value = "2x4.5, 1x2.3"

# Update the selection to reflect where the user ended up after editing
$ws.Range("D6").Select()
